# Generate Report for Handoff
#
# The localization-status report re-sorted the "198a267f..." and
# "8940fb44..." file rows (row 2 <-> row 3) on every sheet (Overview,
# zh-cn, de-de). The hyperlink *addresses* stay anchored to their
# original rId (they are not re-sorted by the report generator), but the
# visible cell text / hyperlink display text does move with the sort.
# Additionally the "f944c77d..." row's Latest Handback DateTime got a
# fresh timestamp on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

function Swap-RowsWithHyperlinks($ws, $row1, $row2, $cols) {
    # Capture the hyperlink addresses for every hyperlink-bearing column,
    # indexed by row, BEFORE we touch any cell values (addresses are
    # keyed off the original row position and must stay put).
    $addr1 = @{}
    $addr2 = @{}
    foreach ($col in $cols) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        if ($c1.Hyperlinks.Count -gt 0) { $addr1[$col] = $c1.Hyperlinks.Item(1).Address }
        if ($c2.Hyperlinks.Count -gt 0) { $addr2[$col] = $c2.Hyperlinks.Item(1).Address }
    }

    # Swap every tracked cell's value between the two rows.
    $allCols = @(1,2,3,4,5,6,7,8,9)
    foreach ($col in $allCols) {
        $v1 = $ws.Cells.Item($row1, $col).Value()
        $v2 = $ws.Cells.Item($row2, $col).Value()
        $ws.Cells.Item($row1, $col).Value = $v2
        $ws.Cells.Item($row2, $col).Value = $v1
    }

    # Rebuild the hyperlinks: same address as before (tied to the row
    # position), but display text now matches the (swapped) cell value.
    $ws.Hyperlinks.Delete()

    $allRels = @()
    for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
        foreach ($col in $cols) {
            if ($r -eq $row1 -and $addr1.ContainsKey($col)) {
                $allRels += ,@($r, $col, $addr1[$col])
            } elseif ($r -eq $row2 -and $addr2.ContainsKey($col)) {
                $allRels += ,@($r, $col, $addr2[$col])
            } elseif ($r -ne $row1 -and $r -ne $row2) {
                $cell = $ws.Cells.Item($r, $col)
                if ($cell.Hyperlinks.Count -gt 0) {
                    $allRels += ,@($r, $col, $cell.Hyperlinks.Item(1).Address)
                }
            }
        }
    }
    foreach ($rel in $allRels) {
        $r = $rel[0]; $col = $rel[1]; $addr = $rel[2]
        $cell = $ws.Cells.Item($r, $col)
        $ws.Hyperlinks.Add($cell, $addr, "", "", $cell.Value())
    }
}

# ---- Overview sheet: swap rows 2 & 3 (column A only has hyperlinks) ----
$wsOverview = $wb.Worksheets.Item("Overview")
Swap-RowsWithHyperlinks $wsOverview 2 3 @(1)

# ---- zh-cn sheet: swap rows 2 & 3 (columns A & C have hyperlinks) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Swap-RowsWithHyperlinks $wsZhCn 2 3 @(1,3)

# ---- de-de sheet: swap rows 2 & 3 (columns A & C have hyperlinks) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
Swap-RowsWithHyperlinks $wsDeDe 2 3 @(1,3)

# ---- f944c77d row (row 5): refreshed Latest Handback DateTime ----
$wsZhCn.Cells.Item(5, 4).Value = "2016-03-03 10:03:07"
$wsDeDe.Cells.Item(5, 4).Value = "2016-03-03 10:03:17"
